$d = $word.ActiveDocument

$replacements = @(
    @("813×3=", "967×2="),
    @("762×4=", "516×2="),
    @("534×4=", "913×5="),
    @("180×8=", "632×7="),
    @("245×2=", "764×6="),
    @("798×5=", "746×7="),
    @("660×2=", "479×2="),
    @("905×4=", "780×8="),
    @("362×7=", "209×4="),
    @("647×8=", "475×8="),
    @("356×2=", "647×9="),
    @("659×8=", "102×7="),
    @("740×2=", "310×3="),
    @("149×7=", "801×7="),
    @("972×6=", "863×2="),
    @("144×4=", "941×4="),
    @("227×5=", "640×6="),
    @("127×7=", "396×9="),
    @("373×6=", "898×7="),
    @("229×4=", "734×8="),
    @("509×8=", "862×4="),
    @("257×2=", "288×2="),
    @("487×6=", "567×5="),
    @("700×2=", "143×9="),
    @("144×5=", "841×9=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
